# Rename the inline picture shapes' display names:
#   - the two Pearson logo images (in the "first page" and "default"
#     footers) from "image1.png" -> "image2.png"
#   - the BTec logo image (in the "first page" header) from
#     "image2.jpg" -> "image1.jpg"
#
# InlineShape does not expose a writable Name in the Word object model,
# so each picture is momentarily converted to a floating Shape (which
# does expose Name), renamed, and converted back to an inline shape so
# the drawing stays wrapped in <wp:inline> exactly as before.

$d = $word.ActiveDocument
$sec = $d.Sections.First

function Rename-InlinePicture($range, $newName) {
    $shape = $range.InlineShapes.Item(1).ConvertToShape()
    $shape.Name = $newName
    $shape.ConvertToInlineShape() | Out-Null
}

# "First page" footer -> footer1.xml (Pearson logo)
Rename-InlinePicture $sec.Footers.Item(2).Range "image2.png"

# "Default" footer -> footer2.xml (Pearson logo)
Rename-InlinePicture $sec.Footers.Item(1).Range "image2.png"

# "First page" header -> header1.xml (BTec logo)
Rename-InlinePicture $sec.Headers.Item(2).Range "image1.jpg"

Write-Output "Renamed inline picture shapes."
